$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.046.39"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "1.650.02"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.98"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.61"
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.265"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").Value = "1.885.16"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "1.658.80"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.79"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "28.050.10"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.96"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +5.56%  "
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.60"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.79"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("D33").Value = "1.448.05"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.893"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.930"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.559"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.44"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  +5.64%  "
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").Value = "1.793.68"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "89.20"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  +0.38%  "
